$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A75").Value = 43796
$ws.Range("A75").NumberFormat = "m/d/yyyy"
